$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 283.92856   # H2 set
$ws.Cells.Item(2, 9).Value = 433.25   # I2 set
$ws.Cells.Item(2, 10).Value = 84.833336   # J2 set
$ws.Cells.Item(2, 11).Value = 433.25   # K2 set
$ws.Cells.Item(2, 12).Value = 84.833336   # L2 set
$ws.Cells.Item(2, 13).Value = -320.25   # M2 set
$ws.Cells.Item(2, 14).Value = -310.833336   # N2 set

$ws.Cells.Item(12, 8).Value = 180.5   # H12 set
$ws.Cells.Item(12, 9).Value = 174.83333   # I12 set
$ws.Cells.Item(12, 10).Value = 197.5   # J12 set
$ws.Cells.Item(12, 11).Value = 174.83333   # K12 set
$ws.Cells.Item(12, 12).Value = 197.5   # L12 set
$ws.Cells.Item(12, 13).Value = -4.833329999999989   # M12 set
$ws.Cells.Item(12, 14).Value = -537.5   # N12 set

$ws.Cells.Item(76, 8).Value = 5977.6   # H76 set
$ws.Cells.Item(76, 9).Value = 5977.6   # I76 set
$ws.Cells.Item(76, 11).Value = 5977.6   # K76 set
$ws.Cells.Item(76, 13).Value = -5662.6   # M76 set

$ws.Cells.Item(79, 8).Value = 5977.6   # H79 set
$ws.Cells.Item(79, 9).Value = 5977.6   # I79 set
$ws.Cells.Item(79, 11).Value = 5977.6   # K79 set
$ws.Cells.Item(79, 13).Value = -4885.6   # M79 set

$ws.Cells.Item(98, 8).Value = 7950.143   # H98 set
$ws.Cells.Item(98, 9).Value = 3214.5715   # I98 set
$ws.Cells.Item(98, 11).Value = 3214.5715   # K98 set
$ws.Cells.Item(98, 13).Value = -1716.5715   # M98 set

$ws.Cells.Item(122, 8).Value = 7950.143   # H122 set
$ws.Cells.Item(122, 9).Value = 3214.5715   # I122 set
$ws.Cells.Item(122, 11).Value = 9643.7145   # K122 set
$ws.Cells.Item(122, 13).Value = -7193.7145   # M122 set

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39, 8).Value = 2000   # H39 set
$ws.Cells.Item(39, 9).Value = 2000   # I39 set
$ws.Cells.Item(39, 10).Value = 0   # J39 set
$ws.Cells.Item(39, 11).Value = 2000   # K39 set
$ws.Cells.Item(39, 12).Value = 0   # L39 set
$ws.Cells.Item(39, 13).Value = -1480   # M39 set
$ws.Cells.Item(39, 14).Value = $null   # N39 delete

$ws.Cells.Item(102, 8).Value = 5002   # H102 set
$ws.Cells.Item(102, 9).Value = 5002   # I102 set
$ws.Cells.Item(102, 11).Value = 5002   # K102 set
$ws.Cells.Item(102, 13).Value = -3380   # M102 set

$ws.Cells.Item(110, 8).Value = 0   # H110 set
$ws.Cells.Item(110, 9).Value = 0   # I110 set
$ws.Cells.Item(110, 11).Value = 0   # K110 set
$ws.Cells.Item(110, 13).Value = $null   # M110 delete

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 7100   # H86 set
$ws.Cells.Item(86, 9).Value = 4000   # I86 set
$ws.Cells.Item(86, 10).Value = 11750   # J86 set
$ws.Cells.Item(86, 11).Value = 4000   # K86 set
$ws.Cells.Item(86, 12).Value = 11750   # L86 set
$ws.Cells.Item(86, 13).Value = -2877   # M86 set
$ws.Cells.Item(86, 14).Value = -13996   # N86 add

$ws.Cells.Item(89, 8).Value = 7100   # H89 set
$ws.Cells.Item(89, 9).Value = 4000   # I89 set
$ws.Cells.Item(89, 10).Value = 11750   # J89 set
$ws.Cells.Item(89, 11).Value = 20000   # K89 set
$ws.Cells.Item(89, 12).Value = 58750   # L89 set
$ws.Cells.Item(89, 13).Value = -14384   # M89 set
$ws.Cells.Item(89, 14).Value = -69982   # N89 add

$ws.Cells.Item(105, 8).Value = 1515.3334   # H105 set
$ws.Cells.Item(105, 9).Value = 1515.3334   # I105 set
$ws.Cells.Item(105, 11).Value = 1515.3334   # K105 set
$ws.Cells.Item(105, 13).Value = 231.6666   # M105 set

$ws.Cells.Item(134, 8).Value = 3166.6667   # H134 set
$ws.Cells.Item(134, 9).Value = 3166.6667   # I134 set
$ws.Cells.Item(134, 11).Value = 9500.000100000001   # K134 set
$ws.Cells.Item(134, 13).Value = -6965.000100000001   # M134 set

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 44.18182   # H7 set
$ws.Cells.Item(7, 9).Value = 48.125   # I7 set
$ws.Cells.Item(7, 10).Value = 33.666668   # J7 set
$ws.Cells.Item(7, 11).Value = 48.125   # K7 set
$ws.Cells.Item(7, 12).Value = 33.666668   # L7 set
$ws.Cells.Item(7, 13).Value = 64.875   # M7 set
$ws.Cells.Item(7, 14).Value = -259.666668   # N7 set

$ws.Cells.Item(22, 8).Value = 228.83333   # H22 set
$ws.Cells.Item(22, 9).Value = 211.8   # I22 set
$ws.Cells.Item(22, 10).Value = 241   # J22 set
$ws.Cells.Item(22, 11).Value = 211.8   # K22 set
$ws.Cells.Item(22, 12).Value = 241   # L22 set
$ws.Cells.Item(22, 13).Value = 138.2   # M22 set
$ws.Cells.Item(22, 14).Value = -941   # N22 set

$ws.Cells.Item(107, 8).Value = 1039.6   # H107 set
$ws.Cells.Item(107, 9).Value = 350   # I107 set
$ws.Cells.Item(107, 10).Value = 1499.3334   # J107 set
$ws.Cells.Item(107, 11).Value = 350   # K107 set
$ws.Cells.Item(107, 12).Value = 1499.3334   # L107 set
$ws.Cells.Item(107, 13).Value = 1570   # M107 set
$ws.Cells.Item(107, 14).Value = -5339.3334   # N107 set

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 200200590   # H4 set
$ws.Cells.Item(4, 9).Value = 250749.75   # I4 set
$ws.Cells.Item(4, 11).Value = 752249.25   # K4 set
$ws.Cells.Item(4, 13).Value = -752137.25   # M4 set

$ws.Cells.Item(106, 8).Value = 0   # H106 set
$ws.Cells.Item(106, 9).Value = 0   # I106 set
$ws.Cells.Item(106, 11).Value = 0   # K106 set
$ws.Cells.Item(106, 13).Value = $null   # M106 delete

$ws.Cells.Item(129, 8).Value = 499   # H129 set
$ws.Cells.Item(129, 9).Value = 499   # I129 set
$ws.Cells.Item(129, 10).Value = 0   # J129 set
$ws.Cells.Item(129, 11).Value = 1497   # K129 set
$ws.Cells.Item(129, 12).Value = 0   # L129 set
$ws.Cells.Item(129, 13).Value = 3503   # M129 add
$ws.Cells.Item(129, 14).Value = $null   # N129 delete

$ws.Cells.Item(131, 8).Value = 3210.5454   # H131 set
$ws.Cells.Item(131, 9).Value = 4343.2   # I131 set
$ws.Cells.Item(131, 10).Value = 2266.6667   # J131 set
$ws.Cells.Item(131, 11).Value = 13029.6   # K131 set
$ws.Cells.Item(131, 12).Value = 6800.000100000001   # L131 set
$ws.Cells.Item(131, 13).Value = -7989.599999999999   # M131 set
$ws.Cells.Item(131, 14).Value = -16880.0001   # N131 set

$ws.Cells.Item(134, 8).Value = 2816   # H134 set
$ws.Cells.Item(134, 9).Value = 1465.75   # I134 set
$ws.Cells.Item(134, 11).Value = 4397.25   # K134 set
$ws.Cells.Item(134, 13).Value = 672.75   # M134 set

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 11638750   # H11 set
$ws.Cells.Item(11, 9).Value = 13835000   # I11 set
$ws.Cells.Item(11, 11).Value = 13835000   # K11 set
$ws.Cells.Item(11, 13).Value = -13834861   # M11 set

$ws.Cells.Item(97, 8).Value = 598.6667   # H97 set
$ws.Cells.Item(97, 10).Value = 598.5   # J97 set
$ws.Cells.Item(97, 12).Value = 598.5   # L97 set
$ws.Cells.Item(97, 14).Value = -1590.5   # N97 set

$ws.Cells.Item(102, 8).Value = 4620.75   # H102 set
$ws.Cells.Item(102, 9).Value = 4994.3335   # I102 set
$ws.Cells.Item(102, 11).Value = 4994.3335   # K102 set
$ws.Cells.Item(102, 13).Value = -3372.3335   # M102 set

$ws.Cells.Item(107, 8).Value = 99   # H107 set
$ws.Cells.Item(107, 9).Value = 99   # I107 set
$ws.Cells.Item(107, 11).Value = 99   # K107 set
$ws.Cells.Item(107, 13).Value = 1821   # M107 set

$ws.Cells.Item(113, 8).Value = 2210.5   # H113 set
$ws.Cells.Item(113, 9).Value = 2210.5   # I113 set
$ws.Cells.Item(113, 10).Value = 0   # J113 set
$ws.Cells.Item(113, 11).Value = 2210.5   # K113 set
$ws.Cells.Item(113, 12).Value = 0   # L113 set
$ws.Cells.Item(113, 13).Value = -40.5   # M113 set
$ws.Cells.Item(113, 14).Value = $null   # N113 delete

$ws.Cells.Item(126, 8).Value = 1912   # H126 set
$ws.Cells.Item(126, 9).Value = 1912   # I126 set
$ws.Cells.Item(126, 11).Value = 5736   # K126 set
$ws.Cells.Item(126, 13).Value = -3266   # M126 set

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 151000   # H19 set
$ws.Cells.Item(19, 9).Value = 151000   # I19 set
$ws.Cells.Item(19, 10).Value = 0   # J19 set
$ws.Cells.Item(19, 11).Value = 151000   # K19 set
$ws.Cells.Item(19, 12).Value = 0   # L19 set
$ws.Cells.Item(19, 13).Value = -150830   # M19 set
$ws.Cells.Item(19, 14).Value = $null   # N19 delete

$ws.Cells.Item(25, 8).Value = 3600   # H25 set
$ws.Cells.Item(25, 9).Value = 3600   # I25 set
$ws.Cells.Item(25, 11).Value = 3600   # K25 set
$ws.Cells.Item(25, 13).Value = -3370   # M25 add

$ws.Cells.Item(29, 8).Value = 10000   # H29 set
$ws.Cells.Item(29, 9).Value = 0   # I29 set
$ws.Cells.Item(29, 10).Value = 10000   # J29 set
$ws.Cells.Item(29, 11).Value = 0   # K29 set
$ws.Cells.Item(29, 12).Value = 10000   # L29 set
$ws.Cells.Item(29, 13).Value = $null   # M29 delete
$ws.Cells.Item(29, 14).Value = -10590   # N29 add

$ws.Cells.Item(100, 8).Value = 2331   # H100 set
$ws.Cells.Item(100, 9).Value = 1997.5   # I100 set
$ws.Cells.Item(100, 10).Value = 2998   # J100 set
$ws.Cells.Item(100, 11).Value = 1997.5   # K100 set
$ws.Cells.Item(100, 12).Value = 2998   # L100 set
$ws.Cells.Item(100, 13).Value = -1456.5   # M100 set
$ws.Cells.Item(100, 14).Value = -4080   # N100 set

$ws.Cells.Item(122, 8).Value = 2877.25   # H122 set
$ws.Cells.Item(122, 9).Value = 2669.8333   # I122 set
$ws.Cells.Item(122, 10).Value = 3499.5   # J122 set
$ws.Cells.Item(122, 11).Value = 8009.499899999999   # K122 set
$ws.Cells.Item(122, 12).Value = 10498.5   # L122 set
$ws.Cells.Item(122, 13).Value = -5559.499899999999   # M122 set
$ws.Cells.Item(122, 14).Value = -15398.5   # N122 set

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 10000   # H2 set
$ws.Cells.Item(2, 9).Value = 10000   # I2 set
$ws.Cells.Item(2, 11).Value = 10000   # K2 set
$ws.Cells.Item(2, 13).Value = -9888   # M2 add

$ws.Cells.Item(15, 8).Value = 9999   # H15 set
$ws.Cells.Item(15, 9).Value = 9999   # I15 set
$ws.Cells.Item(15, 11).Value = 9999   # K15 set
$ws.Cells.Item(15, 13).Value = -9711   # M15 add

$ws.Cells.Item(54, 8).Value = 31500   # H54 set
$ws.Cells.Item(54, 10).Value = 31500   # J54 set
$ws.Cells.Item(54, 12).Value = 31500   # L54 set
$ws.Cells.Item(54, 14).Value = -32540   # N54 add

$ws.Cells.Item(96, 8).Value = 3003   # H96 set
$ws.Cells.Item(96, 9).Value = 3003   # I96 set
$ws.Cells.Item(96, 11).Value = 3003   # K96 set
$ws.Cells.Item(96, 13).Value = -1630   # M96 add
